$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-filled "unique values" column (ShortDivisionName) for the case
# where more than one division shares the same short name - collapse the
# duplicated short names down to "(прочие)" and split "СУП" into its own
# child breakdown rows.

$ws.Range("B5").Value = "(прочие)"
$ws.Range("B10").Value = "(прочие)"
$ws.Range("B12").Value = "СУП+ООД"

# Extend the table with the new detail rows, copying the existing data-row
# formatting down first.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C15").PasteSpecial(-4122)

$ws.Range("A13").Value = "Отдел обеспечения деятельности"
$ws.Range("B13").Value = "СУП+ООД"
$ws.Range("C13").Value = 154

$ws.Range("A14").Value = "Отдел научно-исследовательских разработок"
$ws.Range("B14").Value = "СРБ Отдел НИР"
$ws.Range("C14").Value = 154

$ws.Range("A15").Value = "Отдел поддержки"
$ws.Range("B15").Value = "НСКК"
$ws.Range("C15").Value = 154
